# Genetics.xlsx - "update genetics denied, logout provider"
#
# The test fixture's row 2 holds a single generated case record; the
# case id in A2 (sheet "Sheet1") is refreshed to a newly generated
# "CA-XXXXXXXXXX" token. Several candidate ids were generated/tried
# (matching the batch of new entries appended to the shared-string
# table) before the final one was committed to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "CA-5EDNCJRQ"
$ws.Range("A2").Value = "CA-NN0LBHTF"
$ws.Range("A2").Value = "CA-2555KGWB"
$ws.Range("A2").Value = "CA-QJPKWCFS"
$ws.Range("A2").Value = "CA-J9JRJMAM"
$ws.Range("A2").Value = "CA-3U12KDIS"
$ws.Range("A2").Value = "CA-J7JBFQ38"
$ws.Range("A2").Value = "CA-F6X9FP2J"
$ws.Range("A2").Value = "CA-D3HMI2TW"
$ws.Range("A2").Value = "CA-FYUS3P1A"
$ws.Range("A2").Value = "CA-WYNZ34H6"
$ws.Range("A2").Value = "CA-VX1JXGHE"
$ws.Range("A2").Value = "CA-LQM6KOUH"
$ws.Range("A2").Value = "CA-Z529SCHC"
$ws.Range("A2").Value = "CA-DVJ2WQ9O"
$ws.Range("A2").Value = "CA-2FOE2SHQ"
$ws.Range("A2").Value = "CA-X7ZATUEO"
$ws.Range("A2").Value = "CA-GA0RT8GP"
